$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 6-12: column B availability flips from TRUE to FALSE
for ($r = 6; $r -le 12; $r++) {
    $ws.Cells.Item($r, 2).Value = $false
}

# D30: was boolean TRUE, becomes the literal text string "true" (quote-prefix forces text)
$ws.Cells.Item(30, 4).Value = "'true"
$ws.Cells.Item(30, 4).Style = "Normal"

# New rows 32-34
$ws.Cells.Item(32, 1).Value = 10000030
$ws.Cells.Item(32, 2).Value = $true
$ws.Cells.Item(32, 3).Value = "neues Teil"
$ws.Cells.Item(32, 4).Value = "'"
$ws.Cells.Item(32, 4).Style = "Normal"

$ws.Cells.Item(33, 1).Value = 10000032
$ws.Cells.Item(33, 2).Value = $true
$ws.Cells.Item(33, 3).Value = 200
$ws.Cells.Item(33, 4).Value = "'"
$ws.Cells.Item(33, 4).Style = "Normal"

$ws.Cells.Item(34, 1).Value = 10000031
$ws.Cells.Item(34, 2).Value = $true
$ws.Cells.Item(34, 3).Value = "Teil 281"
$ws.Cells.Item(34, 4).Value = "'"
$ws.Cells.Item(34, 4).Style = "Normal"
